# Add the new sheet "gamesheet_04-03-2024" after the last existing sheet
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "gamesheet_04-03-2024"

# Header row (row 1, columns B:H)
$ws.Range("B1").Value = "Team"
$ws.Range("C1").Value = "Home Score"
$ws.Range("D1").Value = "Away Score"
$ws.Range("E1").Value = "BTTS"
$ws.Range("F1").Value = "First Goal Home"
$ws.Range("G1").Value = "First Goal Away"
$ws.Range("H1").Value = "Nation"

# Game rows (row 2-45): index, matchup, scores, probabilities, nation
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Arsenal v Luton"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.5600000000000001
$ws.Range("F2").Value = 0.67
$ws.Range("G2").Value = 0.33
$ws.Range("H2").Value = "England"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Brentford v Brighton"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.6
$ws.Range("F3").Value = 0.62
$ws.Range("G3").Value = 0.38
$ws.Range("H3").Value = "England"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Manchester City v Aston Villa"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.65
$ws.Range("F4").Value = 0.45
$ws.Range("G4").Value = 0.45
$ws.Range("H4").Value = "England"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Liverpool v Sheffield United"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.6
$ws.Range("F5").Value = 0.6899999999999999
$ws.Range("G5").Value = 0.31
$ws.Range("H5").Value = "England"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Chelsea v Manchester United"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0.75
$ws.Range("F6").Value = 0.12
$ws.Range("G6").Value = 0.75
$ws.Range("H6").Value = "England"
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Crystal Palace v Manchester City"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 0.77
$ws.Range("F7").Value = 0.42
$ws.Range("G7").Value = 0.58
$ws.Range("H7").Value = "England"
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Aston Villa v Brentford"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 0.75
$ws.Range("F8").Value = 0.54
$ws.Range("G8").Value = 0.46
$ws.Range("H8").Value = "England"
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Everton v Burnley"
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = "England"
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Fulham v Newcastle United"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0.51
$ws.Range("F10").Value = 0.6899999999999999
$ws.Range("G10").Value = 0.24
$ws.Range("H10").Value = "England"
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Luton v Bournemouth"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0.43
$ws.Range("F11").Value = 0.45
$ws.Range("G11").Value = 0.55
$ws.Range("H11").Value = "England"
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Wolverhampton Wanderers v West Ham"
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0.62
$ws.Range("F12").Value = 0.75
$ws.Range("G12").Value = 0.25
$ws.Range("H12").Value = "England"
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Brighton v Arsenal"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0.75
$ws.Range("F13").Value = 0.58
$ws.Range("G13").Value = 0.29
$ws.Range("H13").Value = "England"
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Manchester United v Liverpool"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 0.6899999999999999
$ws.Range("F14").Value = 0.58
$ws.Range("G14").Value = 0.42
$ws.Range("H14").Value = "England"
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Sheffield United v Chelsea"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0.61
$ws.Range("F15").Value = 0.33
$ws.Range("G15").Value = 0.57
$ws.Range("H15").Value = "England"
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Tottenham v Nottingham Forest"
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0.83
$ws.Range("F16").Value = 0.54
$ws.Range("G16").Value = 0.46
$ws.Range("H16").Value = "England"
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Lille v Marseille"
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0.39
$ws.Range("F17").Value = 0.78
$ws.Range("G17").Value = 0.23
$ws.Range("H17").Value = "France"
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Lens v Le Havre"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.42
$ws.Range("F18").Value = 0.83
$ws.Range("G18").Value = 0.08
$ws.Range("H18").Value = "France"
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Paris Saint Germain v Clermont Foot"
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 0.59
$ws.Range("F19").Value = 0.76
$ws.Range("G19").Value = 0.17
$ws.Range("H19").Value = "France"
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "Brest v Metz"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0.25
$ws.Range("F20").Value = 0.88
$ws.Range("G20").Value = 0.12
$ws.Range("H20").Value = "France"
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "Montpellier v Lorient"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0.67
$ws.Range("F21").Value = 0.33
$ws.Range("G21").Value = 0.53
$ws.Range("H21").Value = "France"
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "Reims v Nice"
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0.3
$ws.Range("F22").Value = 0.76
$ws.Range("G22").Value = 0.17
$ws.Range("H22").Value = "France"
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "Toulouse v Strasbourg"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0.6
$ws.Range("F23").Value = 0.38
$ws.Range("G23").Value = 0.54
$ws.Range("H23").Value = "France"
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "Monaco v Rennes"
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0.27
$ws.Range("F24").Value = 0.29
$ws.Range("G24").Value = 0.44
$ws.Range("H24").Value = "France"
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "Nantes v Lyon"
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0.6
$ws.Range("F25").Value = 0.3
$ws.Range("G25").Value = 0.7
$ws.Range("H25").Value = "France"
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "Eintracht Frankfurt v Werder Bremen"
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 0.86
$ws.Range("F26").Value = 0.6899999999999999
$ws.Range("G26").Value = 0.31
$ws.Range("H26").Value = "Germany"
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "FC Cologne v Bochum"
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0.75
$ws.Range("F27").Value = 0.62
$ws.Range("G27").Value = 0.25
$ws.Range("H27").Value = "Germany"
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "Freiburg v RasenBallsport Leipzig"
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0.42
$ws.Range("F28").Value = 0.67
$ws.Range("G28").Value = 0.33
$ws.Range("H28").Value = "Germany"
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "Mainz 05 v Darmstadt"
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 0.66
$ws.Range("F29").Value = 0.66
$ws.Range("G29").Value = 0.27
$ws.Range("H29").Value = "Germany"
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "FC Heidenheim v Bayern Munich"
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 0.86
$ws.Range("F30").Value = 0.32
$ws.Range("G30").Value = 0.68
$ws.Range("H30").Value = "Germany"
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "Union Berlin v Bayer Leverkusen"
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 0.5
$ws.Range("F31").Value = 0.5
$ws.Range("G31").Value = 0.38
$ws.Range("H31").Value = "Germany"
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "Borussia Dortmund v VfB Stuttgart"
$ws.Range("C32").Value = 2
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 0.5
$ws.Range("F32").Value = 0.25
$ws.Range("G32").Value = 0.75
$ws.Range("H32").Value = "Germany"
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "Hoffenheim v Augsburg"
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = 0.72
$ws.Range("F33").Value = 0.44
$ws.Range("G33").Value = 0.5600000000000001
$ws.Range("H33").Value = "Germany"
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = "Wolfsburg v Borussia M.Gladbach"
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 0.83
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = "Germany"
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "Salernitana v Sassuolo"
$ws.Range("C35").Value = 2
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 0.55
$ws.Range("F35").Value = 0.62
$ws.Range("G35").Value = 0.38
$ws.Range("H35").Value = "Italy"
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "AC Milan v Lecce"
$ws.Range("C36").Value = 2
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0.49
$ws.Range("F36").Value = 0.84
$ws.Range("G36").Value = 0.16
$ws.Range("H36").Value = "Italy"
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "Roma v Lazio"
$ws.Range("C37").Value = 2
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 0.62
$ws.Range("F37").Value = 0.61
$ws.Range("G37").Value = 0.39
$ws.Range("H37").Value = "Italy"
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = "Empoli v Torino"
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 0.31
$ws.Range("F38").Value = 0.25
$ws.Range("G38").Value = 0.5
$ws.Range("H38").Value = "Italy"
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = "Frosinone v Bologna"
$ws.Range("C39").Value = 1
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 0.73
$ws.Range("F39").Value = 0.42
$ws.Range("G39").Value = 0.42
$ws.Range("H39").Value = "Italy"
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = "Monza v Napoli"
$ws.Range("C40").Value = 1
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0.41
$ws.Range("F40").Value = 0.42
$ws.Range("G40").Value = 0.41
$ws.Range("H40").Value = "Italy"
$ws.Range("A41").Value = 39
$ws.Range("B41").Value = "Cagliari v Atalanta"
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 0.6899999999999999
$ws.Range("F41").Value = 0.39
$ws.Range("G41").Value = 0.46
$ws.Range("H41").Value = "Italy"
$ws.Range("A42").Value = 40
$ws.Range("B42").Value = "Verona v Genoa"
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 1
$ws.Range("E42").Value = 0.57
$ws.Range("F42").Value = 0.64
$ws.Range("G42").Value = 0.29
$ws.Range("H42").Value = "Italy"
$ws.Range("A43").Value = 41
$ws.Range("B43").Value = "Juventus v Fiorentina"
$ws.Range("C43").Value = 2
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 0.55
$ws.Range("F43").Value = 0.65
$ws.Range("G43").Value = 0.3
$ws.Range("H43").Value = "Italy"
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = "Udinese v Inter"
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 2
$ws.Range("E44").Value = 0.21
$ws.Range("F44").Value = 0.07000000000000001
$ws.Range("G44").Value = 0.93
$ws.Range("H44").Value = "Italy"
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "Granada v Valencia"
$ws.Range("C45").Value = 1
$ws.Range("D45").Value = 2
$ws.Range("E45").Value = 0.66
$ws.Range("F45").Value = 0.07000000000000001
$ws.Range("G45").Value = 0.93
$ws.Range("H45").Value = "Spain"

# Formatting: header row + index column (col A) use bold text, thin box border,
# centered horizontally and top-aligned vertically (matches the other gamesheet tabs)
$headerRange = $ws.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$indexRange = $ws.Range("A2:A45")
$indexRange.Font.Bold = $true
$indexRange.Borders.LineStyle = 1
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
